$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove columns I and J entirely (header row 1 and data rows 2-12) ---
$ws.Range("I1:J12").Clear()

# --- Update data rows 2-12, columns B:H with new values ---
$data = @(
    @(0.874, 0.771, -2.991, 1.021, 45.765, 1608, 6386),
    @(0.95, 0.836, -3.067, 1.036, 39.215, 1626, 6455),
    @(0.917, 0.861, -2.98, 1.004, 33.588, 1454, 5766),
    @(0.977, 0.782, -2.961, 1.005, 28.965, 1377, 5458),
    @(1, 0.758, -3.02, 1.038, 24.679, 1097, 4324),
    @(0.986, 0.747, -2.997, 1.039, 20.086, 725, 2841),
    @(1.068, 0.793, -3.084, 1.06, 15.8, 775, 3038),
    @(1.089, 0.846, -3.08, 1.05, 11.749, 771, 3022),
    @(1.08, 0.886, -3.033, 1.027, 8.24, 731, 2862),
    @(1.052, 0.895, -2.969, 1.007, 4.906, 672, 2627),
    @(1.009, 0.964, -2.995, 1.009, 2.243, 661, 2582)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $col = $c + 2   # column B = 2
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}

# --- Add new rows 13-17: copy column-A style from row 12, then set values ---
$ws.Range("A12").Copy($ws.Range("A13:A17"))

$newRows = @(
    @(11, 1.006, 0.979, -2.982, 0.996, 0.109, 557, 2164),
    @(12, 1.004, 0.977, -2.98, 0.996, 0.104, 551, 2138),
    @(13, 1.006, 0.975, -2.98, 0.996, 0.099, 570, 2210),
    @(14, 1.007, 0.974, -2.983, 0.997, 0.095, 607, 2352),
    @(15, 1.008, 0.97, -2.98, 0.997, 0.091, 654, 2537)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $i + 13
    $vals = $newRows[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $col = $c + 1   # column A = 1
        $ws.Cells.Item($row, $col).Value = $vals[$c]
    }
}

Write-Output "Done"
